# Append run: 2026-02-01 06:45 JST
# Updates the "ランサーズ" (Lancers) listing sheet:
#  - refresh the fetch timestamp for the surviving rows
#  - row 2: new listing replacing the old one (title/price/url/score/skills)
#  - row 3: title shortened (price/url/score/skills unchanged)
#  - row 4: new listing replacing the old one
#  - row 5: new listing replacing the old one
#  - row 6 (old "Power Automate..." listing) is dropped entirely
#  - column B widened, column H narrowed back down

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2026-02-01 06:45:36"

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【時給3000~4000円以上/フルリモート】AI駆動開発でのSaaS開発の開発パートナー"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5483313"
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# --- Row 3 (only timestamp + title change) --------------------------------
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "本人顔ベースのリアルタイム顔変換システム開発"

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【PM/フルスタックエンジニア】新規SaaS開発のパートナー募集"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5483306"
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = "◆開発"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "カフェ掲載プラットフォーム「チャヤドコ」開発(要件定義~ベータ版リリース)"
$ws.Range("D5").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5483311"
$ws.Range("G5").Value = 68
$ws.Range("H5").Value = "◆開発"

# --- Drop old row 6 (Power Automate listing) ------------------------------
$ws.Rows.Item(6).Delete()

# --- Hyperlinks: rebuild F2:F5 so the URL rels point at the new targets ---
# (per-item hyperlink deletion isn't reliable on this host, so clear the
#  whole collection and re-add the four that should remain)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5483313")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5483207")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5483306")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5483311")

# --- Column width tweaks ---------------------------------------------------
$ws.Range("B1").ColumnWidth = 47
$ws.Range("H1").ColumnWidth = 12
